$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    8  = -1
    11 = 2
    13 = -1
    15 = 1
    20 = -1
    22 = -6
    23 = 1
    24 = 3
    25 = -2
    32 = 1
    38 = -2
    48 = 1
    61 = 0
    64 = 1
    68 = 2
    70 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
